$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 8) mirroring the existing row 7 layout:
# Date | File name | Lines
$ws.Range("B8").Value2 = $ws.Range("B7").Value2
$ws.Range("B8").NumberFormat = $ws.Range("B7").NumberFormat
$ws.Range("C8").Value2 = "diagram CMD.png"
$ws.Range("D8").Value2 = 100

# Widen column C to fit the new content
$ws.Range("C8").EntireColumn.ColumnWidth = 16

# Move the active selection to D12 (matches the saved cursor position)
$ws.Range("D12").Select()
